# New "vyskladneno" archive entry logged ahead of the existing history:
# a Celý_GB record (GB #7) removed for "Chyba/Spatne zaskladneno", with a
# note about bad data entry. Existing rows 2-8 shift down to 3-9 unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 8
$lastCol = 14   # columns A..N
$numberColumn = 5   # column E (GB_Cislo) stays a real number

# Shift the existing data rows down by one (bottom-up, so we never clobber a
# row before it's been read). Force text columns to Text format before the
# write so date/time-looking strings ("2025-08-02", "14:44:36", ...) land as
# literal text instead of being auto-converted to real Excel dates, and so
# genuinely blank cells stay blank instead of becoming empty strings.
for ($r = $lastRow; $r -ge 2; $r--) {
    for ($col = 1; $col -le $lastCol; $col++) {
        $src = $ws.Cells.Item($r, $col)
        $dst = $ws.Cells.Item($r + 1, $col)
        if ($col -ne $numberColumn) {
            $dst.NumberFormat = "@"
        }
        $dst.Value2 = $src.Value2
    }
}

# Write the new, newest record into the now-empty row 2.
$newRow = @{
    1  = "2025-08-06"
    2  = "10:56:20"
    3  = "cibul"
    4  = "Celý_GB"
    5  = 7
    6  = "KOMPLETNÍ GB #7"
    7  = $null
    8  = $null
    9  = "0 položek"
    10 = $null
    11 = "Chyba/Špatně zaskladněno"
    12 = "........... | blbě zadaná data"
    13 = "2025-07-27"
    14 = "CCCCCCCC"
}

foreach ($col in $newRow.Keys) {
    $cell = $ws.Cells.Item(2, $col)
    if ($col -ne $numberColumn) {
        $cell.NumberFormat = "@"
    }
    $cell.Value2 = $newRow[$col]
}
